$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.51"
$ws.Range("E2").Value = "'0.75%"
$ws.Range("D3").Value = "'35.71"
$ws.Range("E3").Value = "'-4.52%"
$ws.Range("E4").Value = "'1.59%"
$ws.Range("D5").Value = "'0.07858"
$ws.Range("E5").Value = "'0.38%"
$ws.Range("D6").Value = "'2.122"
$ws.Range("E6").Value = "'-3.56%"
$ws.Range("D7").Value = "'7.916"
$ws.Range("E7").Value = "'-1.58%"
$ws.Range("D8").Value = "'0.9181"
$ws.Range("E8").Value = "'0.50%"
$ws.Range("D9").Value = "'0.09763"
$ws.Range("E9").Value = "'0.52%"
$ws.Range("D10").Value = "'0.1863"
$ws.Range("E10").Value = "'-1.43%"
$ws.Range("D11").Value = "'0.08615"
$ws.Range("E11").Value = "'0.41%"
$ws.Range("E12").Value = "'1.02%"
$ws.Range("D13").Value = "'0.09934"
$ws.Range("E13").Value = "'-0.42%"
$ws.Range("D14").Value = "'0.001429"
$ws.Range("E14").Value = "'-3.50%"
$ws.Range("D15").Value = "'0.005679"
$ws.Range("E15").Value = "'0.64%"
$ws.Range("D16").Value = "'3.457"
$ws.Range("E16").Value = "'-0.32%"
$ws.Range("D17").Value = "'4.100"
$ws.Range("E17").Value = "'1.57%"
$ws.Range("D18").Value = "'2.573"
$ws.Range("E18").Value = "'23.64%"
$ws.Range("E19").Value = "'-1.11%"
$ws.Range("D20").Value = "'5.212"
$ws.Range("E20").Value = "'9.38%"
$ws.Range("D21").Value = "'0.1310"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("E22").Value = "'-0.09%"
$ws.Range("E23").Value = "'-1.66%"
$ws.Range("D24").Value = "'0.005055"
$ws.Range("E24").Value = "'5.35%"
$ws.Range("D25").Value = "'0.001236"
$ws.Range("E25").Value = "'0.57%"
$ws.Range("D27").Value = "'0.0004754"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("D39").Value = "'0.01844"
$ws.Range("E39").Value = "'4.72%"
$ws.Range("D40").Value = "'0.04726"
$ws.Range("E40").Value = "'-0.11%"
$ws.Range("D41").Value = "'0.007480"
$ws.Range("E41").Value = "'-7.29%"
$ws.Range("E42").Value = "'0.35%"
$ws.Range("D43").Value = "'0.007753"
$ws.Range("E43").Value = "'1.73%"
$ws.Range("D44").Value = "'0.002242"
$ws.Range("E44").Value = "'2.87%"
$ws.Range("D45").Value = "'0.01101"
$ws.Range("E45").Value = "'5.60%"
$ws.Range("D46").Value = "'0.00006321"
$ws.Range("E46").Value = "'4.49%"
$ws.Range("D48").Value = "'0.0005806"
$ws.Range("E48").Value = "'0.09%"
$ws.Range("D49").Value = "'47.33"
$ws.Range("E49").Value = "'615.70%"
$ws.Range("E50").Value = "'-25.59%"
$ws.Range("D51").Value = "'0.00002102"
